$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; existing rows 28:90 shift down to 29:91
$ws.Rows(28).Insert()

# Populate the new row 28 with data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R carry over the same values as the (now shifted)
# row below it (old row 28, now row 29); D,J,K,L,M,P hold the new reported values.
$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 45044
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 100112010
$ws.Range("G28").Value = "Achicoria"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 10000
$ws.Range("N28").Value = "$/caja 18 unidades"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 556
$ws.Range("Q28").Value = 18
$ws.Range("R28").Value = "Hortaliza"
